$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.404.74"
$ws.Range("E2").Value = "  -0.06%  "

$ws.Range("D3").Value = "2.069.38"
$ws.Range("E3").Value = "  +0.04%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'235.31"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").Value = "'0.628"
$ws.Range("E6").Value = "  +2.20%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'57.26"
$ws.Range("E8").Value = "  -2.16%  "

$ws.Range("E9").Value = "  +2.83%  "

$ws.Range("D10").Value = "'0.0774"
$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("E11").Value = "  +0.77%  "

$ws.Range("D12").Value = "2.370.35"
$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").Value = "'14.42"
$ws.Range("E13").Value = "  -1.51%  "

$ws.Range("D14").Value = "'20.74"
$ws.Range("E14").Value = "  -1.12%  "

$ws.Range("D15").Value = "'0.776"

$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").Value = "2.067.94"
$ws.Range("E17").Value = "  -0.24%  "

$ws.Range("D18").Value = "37.309.01"
$ws.Range("E18").Value = "  -0.84%  "

$ws.Range("D19").Value = "'6.16"
$ws.Range("E19").Value = "  -0.67%  "

$ws.Range("D20").Value = "'69.56"
$ws.Range("E20").Value = "  +0.70%  "

$ws.Range("D21").Value = "0.0₃0818"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").Value = "'226.87"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("E24").Value = "  +1.63%  "

$ws.Range("D25").Value = "'2.40"
$ws.Range("E25").Value = "  -2.37%  "

$ws.Range("D26").Value = "'167.12"
$ws.Range("E26").Value = "  +1.46%  "

$ws.Range("E27").Value = "  +0.33%  "

$ws.Range("E28").Value = "  -7.06%  "

$ws.Range("D29").Value = "'0.128"
$ws.Range("E29").Value = "  +0.99%  "

$ws.Range("D30").Value = "'19.10"
$ws.Range("E30").Value = "  -0.77%  "

$ws.Range("E31").Value = "  -0.48%  "

$ws.Range("D32").Value = "'4.55"
$ws.Range("E32").Value = "  +0.45%  "

$ws.Range("E33").Value = "  -1.11%  "

$ws.Range("D34").Value = "'4.53"
$ws.Range("E34").Value = "  +0.47%  "

$ws.Range("E35").Value = "  -3.35%  "

$ws.Range("E36").Value = "  +0.46%  "

$ws.Range("D37").Value = "'3.34"
$ws.Range("E37").Value = "  -2.50%  "

$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").Value = "'5.61"
$ws.Range("E39").Value = "  -4.96%  "

$ws.Range("E40").Value = "  -0.28%  "

$ws.Range("D41").Value = "'0.0957"
$ws.Range("E41").Value = "  -2.78%  "

$ws.Range("D42").Value = "1.490.48"
$ws.Range("E42").Value = "  +0.79%  "

$ws.Range("D43").Value = "'97.67"
$ws.Range("E43").Value = "  +0.74%  "

$ws.Range("E44").Value = "  +0.81%  "

$ws.Range("E45").Value = "  -1.13%  "

$ws.Range("D46").Value = "'4.15"
$ws.Range("E46").Value = "  -7.68%  "

$ws.Range("D47").Value = "'1.03"
$ws.Range("E47").Value = "  -0.81%  "

$ws.Range("E48").Value = "  -1.42%  "

$ws.Range("D49").Value = "'15.09"
$ws.Range("E49").Value = "  -5.45%  "

$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("D51").Value = "'47.34"
$ws.Range("E51").Value = "  +5.53%  "
